$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14, pushing existing rows 14-48 down to 15-49.
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new weekly record.
$ws.Cells.Item(14, 1).Value = 5
$ws.Cells.Item(14, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(14, 3).Value = 'Maule'
$ws.Cells.Item(14, 4).Value = 44498
$ws.Cells.Item(14, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(14, 5).Value = 7
$ws.Cells.Item(14, 6).Value = 'Fruta'
$ws.Cells.Item(14, 7).Value = 100107
$ws.Cells.Item(14, 8).Value = 'Otros'
$ws.Cells.Item(14, 9).Value = 100107002
$ws.Cells.Item(14, 10).Value = 'Chirimoya'
$ws.Cells.Item(14, 11).Value = 'Cultivar IV Región'
$ws.Cells.Item(14, 12).Value = 'Primera'
$ws.Cells.Item(14, 13).Value = 250
$ws.Cells.Item(14, 14).Value = 22000
$ws.Cells.Item(14, 15).Value = 23000
$ws.Cells.Item(14, 16).Value = 22600
$ws.Cells.Item(14, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(14, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(14, 19).Value = 2260
$ws.Cells.Item(14, 20).Value = 10
